# Weekly update: insert a new price record at row 376, pushing the
# existing rows 376-431 down to 377-432 (dimension grows from
# A1:R431 to A1:R432), and populate the newly inserted row with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 376; this shifts every
# row from 376 through 431 down by one (to 377-432) and carries the
# existing formatting (e.g. the date number format on column D) along.
$ws.Rows.Item(376).Insert()

# Populate the newly-inserted row 376 with this week's record.
$ws.Range("A376").Value = 4
$ws.Range("B376").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C376").Value = 'Los Lagos'
$ws.Range("D376").Value = 45142
$ws.Range("E376").Value = 10
$ws.Range("F376").Value = 100112044
$ws.Range("G376").Value = 'Perejil'
$ws.Range("H376").Value = 'Sin especificar'
$ws.Range("I376").Value = 'Primera'
$ws.Range("J376").Value = 180
$ws.Range("K376").Value = 6000
$ws.Range("L376").Value = 6000
$ws.Range("M376").Value = 6000
$ws.Range("N376").Value = '$/docena de atados (3 kilos)'
$ws.Range("O376").Value = 'Región Metropolitana'
$ws.Range("P376").Value = 2000
$ws.Range("Q376").Value = 3
$ws.Range("R376").Value = 'Hortaliza'
